$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2410"
$ws.Range("E17").Value = "2411"
$ws.Range("E18").Value = "2412"
$ws.Range("E19").Value = "2501"
$ws.Range("E20").Value = "2502"
$ws.Range("E21").Value = "2503"
$ws.Range("E22").Value = "2504"

$ws.Range("F16").Value = 52000
$ws.Range("F22").Value = 39866
